$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, whether the cell must be forced to
# Text format first so Excel does not coerce numeric-looking strings (e.g.
# "324.45" or "2.040") into floating point numbers and lose exact formatting.
$updates = @(
    @{ Cell = 'D2'; Value = '27.638.78'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -3.37%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.861.48'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -4.69%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  -0.80%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '324.45'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +0.92%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '1.008'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -0.56%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.4498'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  -5.48%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.3863'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -4.16%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '49.27'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -8.40%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.08026'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -4.49%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '1.019'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -3.75%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '21.49'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -2.91%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '1.902.28'; ForceText = $false }
    @{ Cell = 'E13'; Value = '  -1.98%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '7.191'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -5.28%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '5.882'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  -4.99%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '1.013'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  -0.32%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '86.63'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  -2.63%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '0.00001034'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -3.65%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '0.06563'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -1.10%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '17.11'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -8.17%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '1.008'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  -0.46%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '5.512'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -5.15%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '27.613.47'; ForceText = $false }
    @{ Cell = 'E23'; Value = '  -3.59%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '10.87'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -5.47%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '2.312'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  +0.59%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '2.114.15'; ForceText = $false }
    @{ Cell = 'E26'; Value = '  -2.69%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '151.69'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -1.62%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '19.49'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -3.12%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '5.542'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -6.22%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '2.040'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -5.17%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '120.89'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -2.11%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '0.09407'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -1.69%  '; ForceText = $false }
    @{ Cell = 'B33'; Value = 'ARBITRUM'; ForceText = $false }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText = $false }
    @{ Cell = 'D33'; Value = '1.460'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +2.02%  '; ForceText = $false }
    @{ Cell = 'B34'; Value = 'ImmutableX'; ForceText = $false }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; ForceText = $false }
    @{ Cell = 'D34'; Value = '0.9282'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -6.84%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '3.639'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.85%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '5.302'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -4.93%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '1.230'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -2.13%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.02234'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  -4.05%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '0.06002'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -3.57%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '8.409'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  -3.45%  '; ForceText = $false }
    @{ Cell = 'B41'; Value = 'Frax'; ForceText = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; ForceText = $false }
    @{ Cell = 'D41'; Value = '1.007'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -0.52%  '; ForceText = $false }
    @{ Cell = 'B42'; Value = 'TheSandbox'; ForceText = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.5950'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -4.30%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '10.34'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -6.54%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '0.1854'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -3.38%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '1.281'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -3.64%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '12.52'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -3.09%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '0.5656'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -4.64%  '; ForceText = $false }
    @{ Cell = 'B48'; Value = 'PancakeSwap'; ForceText = $false }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; ForceText = $false }
    @{ Cell = 'D48'; Value = '3.417'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +0.03%  '; ForceText = $false }
    @{ Cell = 'B49'; Value = 'NEARProtocol'; ForceText = $false }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; ForceText = $false }
    @{ Cell = 'D49'; Value = '1.934'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -6.50%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.06878'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  +0.80%  '; ForceText = $false }
    @{ Cell = 'B51'; Value = 'PaxosStandard'; ForceText = $false }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'; ForceText = $false }
    @{ Cell = 'D51'; Value = '1.008'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -0.66%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $range.NumberFormat = '@'
    }
    $range.Value = $u.Value
}
